$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 203.7816646666667
$ws.Range("H2").Value = 611.344994
$ws.Range("I2").Value = 0.6667327591988204
$ws.Range("J2").Value = 0.6667327591988205
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 409.6166503333334
$ws.Range("N2").Value = 1228.849951
$ws.Range("O2").Value = 0.6234125531262766
$ws.Range("P2").Value = 0.6234125531262766
$ws.Range("Q2").Value = 83472.3628801106
$ws.Range("R2").Value = 751251.2659209954
$ws.Range("S2").Value = 0.4156495716650636
$ws.Range("T2").Value = 0.4156495716650637

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 203.7816646666667
$ws.Range("H3").Value = 611.344994
$ws.Range("I3").Value = 0.6667327591988204
$ws.Range("J3").Value = 0.6667327591988205
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 56.495384
$ws.Range("N3").Value = 169.486152
$ws.Range("O3").Value = 0.08598266586728959
$ws.Range("P3").Value = 0.08598266586728959
$ws.Range("Q3").Value = 11512.72339750257
$ws.Range("R3").Value = 103614.5105775231
$ws.Range("S3").Value = 0.05732746005696822
$ws.Range("T3").Value = 0.05732746005696823

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 203.7816646666667
$ws.Range("H4").Value = 611.344994
$ws.Range("I4").Value = 0.6667327591988204
$ws.Range("J4").Value = 0.6667327591988205
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 190.9434713333333
$ws.Range("N4").Value = 572.830414
$ws.Range("O4").Value = 0.2906047810064339
$ws.Range("P4").Value = 0.2906047810064338
$ws.Range("Q4").Value = 38910.77844553861
$ws.Range("R4").Value = 350197.0060098476
$ws.Range("S4").Value = 0.1937557274767886
$ws.Range("T4").Value = 0.1937557274767886

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("H5").Value = 189.421768
$ws.Range("I5").Value = 0.2065833519051582
$ws.Range("J5").Value = 0.2065833519051582
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 409.6166503333334
$ws.Range("N5").Value = 1228.849951
$ws.Range("O5").Value = 0.6234125531262766
$ws.Range("P5").Value = 0.6234125531262766
$ws.Range("Q5").Value = 25863.4367027926
$ws.Range("R5").Value = 232770.9303251334
$ws.Range("S5").Value = 0.1287866548445787
$ws.Range("T5").Value = 0.1287866548445787

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 63.14058933333333
$ws.Range("H6").Value = 189.421768
$ws.Range("I6").Value = 0.2065833519051582
$ws.Range("J6").Value = 0.2065833519051582
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 56.495384
$ws.Range("N6").Value = 169.486152
$ws.Range("O6").Value = 0.08598266586728959
$ws.Range("P6").Value = 0.08598266586728959
$ws.Range("Q6").Value = 3567.151840372971
$ws.Range("R6").Value = 32104.36656335674
$ws.Range("S6").Value = 0.01776258732060592
$ws.Range("T6").Value = 0.01776258732060592

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 63.14058933333333
$ws.Range("H7").Value = 189.421768
$ws.Range("I7").Value = 0.2065833519051582
$ws.Range("J7").Value = 0.2065833519051582
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 190.9434713333333
$ws.Range("N7").Value = 572.830414
$ws.Range("O7").Value = 0.2906047810064339
$ws.Range("P7").Value = 0.2906047810064338
$ws.Range("Q7").Value = 12056.28330933911
$ws.Range("R7").Value = 108506.549784052
$ws.Range("S7").Value = 0.06003410973997356
$ws.Range("T7").Value = 0.06003410973997355

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 38.719942
$ws.Range("H8").Value = 116.159826
$ws.Range("I8").Value = 0.1266838888960214
$ws.Range("J8").Value = 0.1266838888960214
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 409.6166503333334
$ws.Range("N8").Value = 1228.849951
$ws.Range("O8").Value = 0.6234125531262766
$ws.Range("P8").Value = 0.6234125531262766
$ws.Range("Q8").Value = 15860.33294314095
$ws.Range("R8").Value = 142742.9964882685
$ws.Range("S8").Value = 0.07897632661663426
$ws.Range("T8").Value = 0.07897632661663428

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 38.719942
$ws.Range("H9").Value = 116.159826
$ws.Range("I9").Value = 0.1266838888960214
$ws.Range("J9").Value = 0.1266838888960214
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 56.495384
$ws.Range("N9").Value = 169.486152
$ws.Range("O9").Value = 0.08598266586728959
$ws.Range("P9").Value = 0.08598266586728959
$ws.Range("Q9").Value = 2187.497991747728
$ws.Range("R9").Value = 19687.48192572955
$ws.Range("S9").Value = 0.01089261848971545
$ws.Range("T9").Value = 0.01089261848971545

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 38.719942
$ws.Range("H10").Value = 116.159826
$ws.Range("I10").Value = 0.1266838888960214
$ws.Range("J10").Value = 0.1266838888960214
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 190.9434713333333
$ws.Range("N10").Value = 572.830414
$ws.Range("O10").Value = 0.2906047810064339
$ws.Range("P10").Value = 0.2906047810064338
$ws.Range("Q10").Value = 7393.320135305329
$ws.Range("R10").Value = 66539.88121774797
$ws.Range("S10").Value = 0.0368149437896717
$ws.Range("T10").Value = 0.0368149437896717
